# Update database and change read_price algorithm:
# Shift the 5-year rolling window forward by one year (drop 1396/12, add 1401/12),
# sliding existing E:H values left into D:G... i.e. E<-F, F<-G, G<-H, H<-I, and I gets new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows: year labels (row 8 and row 24) ---
$ws.Range("E8").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Data rows: shift existing F:I values left into E:H, then put new value into I ---
$dataRows = @{
    10 = 1277218
    11 = 0
    12 = 0
    13 = 514316
    14 = 0
    15 = 549
    16 = 42540
    17 = 607849
    18 = 0
    19 = 813757
    20 = 3256229
    26 = 435
    27 = 311
}

foreach ($r in $dataRows.Keys) {
    $fVal = $ws.Cells.Item($r, 6).Value2
    $gVal = $ws.Cells.Item($r, 7).Value2
    $hVal = $ws.Cells.Item($r, 8).Value2
    $iVal = $ws.Cells.Item($r, 9).Value2

    $ws.Cells.Item($r, 5).Value2 = $fVal
    $ws.Cells.Item($r, 6).Value2 = $gVal
    $ws.Cells.Item($r, 7).Value2 = $hVal
    $ws.Cells.Item($r, 8).Value2 = $iVal
    $ws.Cells.Item($r, 9).Value2 = $dataRows[$r]
}
